$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-8 from 45208 to 45212
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45212
}
